# Add "prolificid" column (offer gender in binary => adds anonymized participant id column)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts name/gender/matrices/race/mat_rank one column to the right
$ws.Columns.Item(3).Insert()

# Give the new header cell (C1) the same look (bold/centered/bordered) as the other header cells
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 3).Value = "prolificid"

# Final table (rows re-ranked after refreshed "matrices" scores, each row tagged with its prolific id)
$rows = @(
    @{B=3;  C='60bd88b8fc436774352f53b9'; D='Annes';     E='female'; F=13.44015278694422; G='Asian';                      H=1},
    @{B=2;  C='5c5882fc5bfe7600011197cb'; D='Colleen';   E='female'; F=13.26073009200145; G='White';                      H=2},
    @{B=21; C='5c0e89c6c323400001e6c4a5'; D='Bri';       E='female'; F=8.354532088468041; G='Black or African American';  H=3},
    @{B=19; C='60b45e9961dd412bfb6780f8'; D='Jewel';     E='female'; F=8.098102820749885; G='Black or African American';  H=4},
    @{B=22; C='608b14a312c099ac00b721b6'; D='Khushi';    E='female'; F=8.097338185867613; G='Asian';                      H=5},
    @{B=33; C='60cb36ee9f58331a33cf5506'; D='Shaniek';   E='female'; F=5.483549169353528; G='Black or African American';  H=6},
    @{B=32; C='6036f9b3b1842f8b659b18c7'; D='Kellie';    E='female'; F=5.36924149709817;  G='White';                      H=7},
    @{B=30; C='60d5775a99b502eec8cf56b4'; D='Shadaisia'; E='female'; F=5.049422940202584; G='Black or African American';  H=8},
    @{B=34; C='5e96194b0a9fe909389e9f7b'; D='Tina';      E='female'; F=4.105146646021751; G='White';                      H=9},
    @{B=35; C='6077db0613ce87b4a62a78f9'; D='Lori';      E='female'; F=4.050914323979571; G='White';                      H=10},
    @{B=41; C='60bfcf5805c5ae12a546f9f3'; D='Giana';     E='female'; F=2.435214467044919; G='White';                      H=11},
    @{B=44; C='60c0e5899d387663c07eb3a4'; D='Nansi';     E='female'; F=1.191712437135525; G='Asian';                      H=12}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $r = $r + 1
}
